# ValueSet-KLInformationGatheringTypeCodes.xlsx : bump term version to 1.1.0
#
# Semantic changes performed (per the commit's xml diff):
#   1. Metadata!B3 ("Version" row)  1.0.0 -> 1.1.0
#   2. Metadata!B8 ("Date" row)     2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00
#   3. Re-assert wrap/vertical-top alignment (applyAlignment) on the header
#      row style and the body-row style on both worksheets, so the saved
#      cellXfs entries carry applyAlignment="true" alongside the existing
#      vertical="top" wrapText="true" <alignment>.

$wb = $excel.ActiveWorkbook

# ---- 1 & 2: update the Version and Date values on the Metadata sheet ----
$ws1 = $wb.Worksheets.Item("Metadata")

$versionLabel = $ws1.Cells.Find("Version")
$versionCell = $versionLabel.Offset(0, 1)
$versionCell.Value = "1.1.0"

$dateLabel = $ws1.Cells.Find("Date")
$dateCell = $dateLabel.Offset(0, 1)
$dateCell.Value = "2023-07-10T23:08:03+02:00"

# ---- 3: re-apply alignment/wrap so applyAlignment is written on save ----
# Metadata sheet: header row (A1:B1) + body rows (A2:B14)
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A1:B1").VerticalAlignment = -4160   # xlTop

$ws1.Range("A2:B14").WrapText = $true
$ws1.Range("A2:B14").VerticalAlignment = -4160  # xlTop

# "Include from CareSocialCodes" sheet: header row + the two populated
# body blocks (row 2 spans A:C, rows 3-4 only span A:B in the source file)
$ws2 = $wb.Worksheets.Item("Include from CareSocialCodes")

$ws2.Range("A1:C1").WrapText = $true
$ws2.Range("A1:C1").VerticalAlignment = -4160

$ws2.Range("A2:C2").WrapText = $true
$ws2.Range("A2:C2").VerticalAlignment = -4160

$ws2.Range("A3:B4").WrapText = $true
$ws2.Range("A3:B4").VerticalAlignment = -4160
